$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row at position 4 for "course_resultat", shifting cheval_nom etc. down.
$ws.Rows(4).Insert(-4121)

# Give the new row 4 the same "blue" formatting as row 3 (an existing course_* row).
$ws.Range("A3:E3").Copy()
$ws.Range("A4:E4").PasteSpecial(-4122)

$ws.Range("A4").Value = "course_resultat"
$ws.Range("B4").Value = 1

# 2. Remove the old "cheval_position_arrivee" row (now row 7, after the insert shifted it down).
$ws.Rows(7).Delete()

# 3. Recolor rows 10 and 11 ("pari_cheval_position"/"type_id", about to become the new
#    "categorie_id"/"categorie_libelle" lookup table) with the existing "pink" style used by
#    row 12 ("type_libelle"), before that template row is removed.
$ws.Range("A12:E12").Copy()
$ws.Range("A10:E10").PasteSpecial(-4122)
$ws.Range("A12:E12").Copy()
$ws.Range("A11:E11").PasteSpecial(-4122)

# 4. Rename pari_gain -> pari_ordre (row 9, value/position unchanged).
$ws.Range("A9").Value = "pari_ordre"

# 5. Rename pari_cheval_position -> categorie_id (row 10, value/position unchanged).
$ws.Range("A10").Value = "categorie_id"

# 6. Rename type_id -> categorie_libelle (row 11), moving its "1" marker from column D to column E.
$ws.Range("A11").Value = "categorie_libelle"
$ws.Range("D11").ClearContents()
$ws.Range("E11").Value = 1

# 7. Remove the old trailing "type_libelle" row (now row 12).
$ws.Rows(12).Delete()

# 8. Update the header row: E1 now references categorie_id instead of type_id.
$ws.Range("E1").Value = "categorie_id"

# 9. Column E's best-fit width grows slightly to fit "categorie_id".
$ws.Columns("E:E").ColumnWidth = 11.166666666666666

# 10. Leave the selection where the author last left it when saving.
[void]$ws.Range("E17").Select()
